$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Primer-Probe Interactions")

$ws.Range("E6").Value = "GTCCTCACTGCTTTTCAACCCTAT"
$ws.Range("F6").Value = "TGCAACCTGTCCTTTATTTTTCC"
$ws.Range("D6").Value = "CCCAGAGCTCTCTAC"

$ws.Activate()
$ws.Range("F16").Select()
